# Applies the cryptocurrency price/volume update described in the commit:
# "Updated cryptos list on Sun Jul 30 22:35:35 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.191.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.859.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7031'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3106'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07999'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.858.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.170'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6950'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.354'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.159.72'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008279'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '250.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.115.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('E21').Value = '  -1.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.500'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1548'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.960'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.496'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.251'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05244'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.878'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7413'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.17%  '
$ws.Range('E36').Value = '  -2.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.710'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01861'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.244.70'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.733'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.234'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8949'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000128'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.015.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.777'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.401'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('E51').Value = '  -2.60%  '
